$wb = $excel.ActiveWorkbook

# --- TC09 (sheet10.xml): password "Test" -> "WSCAdmin"; selection A1:E2 -> A2:B2
$wsTC09 = $wb.Worksheets.Item("TC09")
$wsTC09.Range("B2").Value = "WSCAdmin"
$wsTC09.Range("A2:B2").Select()

# --- TC08 (sheet8.xml): password "Test" -> "WSCAdmin"; selection B2 -> A2:B2
$wsTC08 = $wb.Worksheets.Item("TC08")
$wsTC08.Range("B2").Value = "WSCAdmin"
$wsTC08.Range("A2:B2").Select()

# --- TC10 (sheet11.xml): no direct cell/selection edits (only shared-string reindex, handled automatically)

# --- TC11 (sheet12.xml): was the active tab; selection F2 -> A2:B2, loses tabSelected
$wsTC11 = $wb.Worksheets.Item("TC11")
$wsTC11.Range("A2:B2").Select()

# --- TC13 (sheet14.xml): no selection change, only shared-string reindex (handled automatically)

# --- TC14 (sheet15.xml): password "Test" -> "WSCAdmin"; selection B2 -> A41
$wsTC14 = $wb.Worksheets.Item("TC14")
$wsTC14.Range("B2").Value = "WSCAdmin"
$wsTC14.Range("A41").Select()

# --- TC02 (sheet2.xml): becomes the new active tab; RunAs label column values collapse
#     to duplicate the RunAs column to their right-hand neighbour's value; column D widens.
$wsTC02 = $wb.Worksheets.Item("TC02")
$wsTC02.Range("E2").Value = "UTMG"
$wsTC02.Range("E3").Value = "CWC"
$wsTC02.Range("E4").Value = "Jordan Valley"
$wsTC02.Columns.Item(4).ColumnWidth = 13.5

# --- TC06 (sheet6.xml): drop frozen/scrolled topLeftCell, selection K2 -> A2
$wsTC06 = $wb.Worksheets.Item("TC06")
$wsTC06.Range("A2").Select()

# Activate TC02 last so it becomes the workbook's active tab (activeTab goes 11 -> 1)
# and picks up tabSelected="1" on its sheetView, matching the commit's view-state fix.
$wsTC02.Activate()
$wsTC02.Range("A2:B2").Select()

$wb.Save()
